$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $val) {
    $c = $ws.Range($cellAddr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.NumberFormat = "General"
    $c.Style = "Normal"
}

Set-TextValue "D2" "304.12"
Set-TextValue "E2" "1.91%"
Set-TextValue "D3" "31.65"
Set-TextValue "E3" "-0.09%"
Set-TextValue "D4" "5.166"
Set-TextValue "E4" "2.12%"
Set-TextValue "D5" "0.07496"
Set-TextValue "E5" "-0.32%"
Set-TextValue "D6" "2.391"
Set-TextValue "E6" "39.41%"
Set-TextValue "D7" "8.010"
Set-TextValue "E7" "2.80%"
Set-TextValue "D8" "3.865"
Set-TextValue "E8" "1.83%"
Set-TextValue "D9" "0.9152"
Set-TextValue "E9" "-1.19%"
Set-TextValue "E10" "1.92%"
Set-TextValue "D11" "0.07702"
Set-TextValue "E11" "3.31%"
Set-TextValue "D12" "0.08134"
Set-TextValue "E12" "1.99%"
Set-TextValue "D13" "0.03039"
Set-TextValue "E13" "-0.22%"
Set-TextValue "D14" "0.09934"
Set-TextValue "E14" "0.45%"
Set-TextValue "E15" "0.46%"
Set-TextValue "D16" "0.006090"
Set-TextValue "E16" "-3.98%"
Set-TextValue "D17" "3.502"
Set-TextValue "E17" "1.47%"
Set-TextValue "D18" "2.239"
Set-TextValue "E18" "0.96%"
Set-TextValue "D19" "0.3263"
Set-TextValue "D20" "0.1328"
Set-TextValue "E20" "-0.56%"
Set-TextValue "D21" "4.658"
Set-TextValue "E21" "2.26%"
Set-TextValue "D22" "0.04594"
Set-TextValue "E22" "-1.46%"
Set-TextValue "E24" "3.86%"
Set-TextValue "D25" "0.004536"
Set-TextValue "E25" "2.43%"
Set-TextValue "E26" "-7.25%"
Set-TextValue "D27" "0.0002737"
Set-TextValue "E27" "48.50%"
Set-TextValue "D39" "0.01753"
Set-TextValue "E39" "5.11%"
Set-TextValue "D40" "0.04540"
Set-TextValue "E40" "-0.24%"
Set-TextValue "D41" "0.007408"
Set-TextValue "E41" "5.59%"
Set-TextValue "D42" "0.1364"
Set-TextValue "E42" "2.82%"
Set-TextValue "D43" "0.002158"
Set-TextValue "E43" "4.73%"
Set-TextValue "E44" "-14.90%"
Set-TextValue "D45" "0.00006543"
Set-TextValue "E45" "7.46%"
Set-TextValue "E46" "15.31%"
